$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.106
$ws.Range("D3").Value = -0.106
$ws.Range("G2").Value = -37
$ws.Range("G3").Value = -37
$ws.Range("H2").Value = -37
$ws.Range("H3").Value = -37
$ws.Range("I2").Value = -36.08333333333334
$ws.Range("I3").Value = -36.08333333333334
$ws.Range("J2").Value = -36.08333333333334
$ws.Range("J3").Value = -36.08333333333334
$ws.Range("K2").Value = -0.605
$ws.Range("K3").Value = -0.605
$ws.Range("L2").Value = -50.41666666666666
$ws.Range("L3").Value = -50.41666666666666
$ws.Range("U2").Value = 0.108
$ws.Range("U3").Value = 0.108
$ws.Range("V2").Value = 0.001463414634146341
$ws.Range("V3").Value = 0.001463414634146341
$ws.Range("W2").Value = -0.136568848758465
$ws.Range("W3").Value = -0.136568848758465
$ws.Range("X2").Value = 0.05387377572522387
$ws.Range("X3").Value = 0.05387377572522387
$ws.Range("Y2").Value = -0.1904426244836889
$ws.Range("Y3").Value = -0.1904426244836889
$ws.Range("Z2").Value = 0.001628222523744912
$ws.Range("Z3").Value = 0.001628222523744912
$ws.Range("AA2").Value = -0.0587516960651289
$ws.Range("AA3").Value = -0.0587516960651289
$ws.Range("AB2").Value = 0.0533858702594333
$ws.Range("AB3").Value = 0.0533858702594333
$ws.Range("AC2").Value = -0.1121375663245622
$ws.Range("AC3").Value = -0.1121375663245622
$ws.Range("AD2").Value = 3.75
$ws.Range("AD3").Value = 3.75
$ws.Range("AF2").Value = 3.75
$ws.Range("AF3").Value = 3.75
$ws.Range("AG2").Value = 3.642
$ws.Range("AG3").Value = 3.642
$ws.Range("AH2").Value = 0.04835589941972921
$ws.Range("AH3").Value = 0.04835589941972921
$ws.Range("AI2").Value = 0.5047106325706595
$ws.Range("AI3").Value = 0.5047106325706595
$ws.Range("AJ2").Value = 0.04702874409235299
$ws.Range("AJ3").Value = 0.04702874409235299
$ws.Range("AK2").Value = 0.4974050805790767
$ws.Range("AK3").Value = 0.4974050805790767
$ws.Range("AL2").Value = 0.201
$ws.Range("AL3").Value = 0.201
$ws.Range("AM2").Value = 0.189
$ws.Range("AM3").Value = 0.189
$ws.Range("AN2").Value = -8.761682242990654
$ws.Range("AN3").Value = -8.761682242990654
$ws.Range("AO2").Value = -2.154228855721393
$ws.Range("AO3").Value = -2.154228855721393
$ws.Range("AP2").Value = -8.509345794392523
$ws.Range("AP3").Value = -8.509345794392523
$ws.Range("AQ2").Value = -2.291005291005291
$ws.Range("AQ3").Value = -2.291005291005291
